# Natmi following Dr Hou advice
# Rebuild the LR-pair result table (rows 2-17) as a full 4x4 cross-join
# of sending/target clusters (ECs, FAPs, M2, sCs) for ligand Postn / receptor Ptk7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array: Sending cluster, Ligand, Receptor, Target cluster, then the
# 16 numeric metric columns (E..T), all as strings (cast to double on write so
# scientific-notation values parse correctly).
$data = @(
    @("ECs", "Postn", "Ptk7", "ECs", "3", "1", "7.018436333333334", "21.055309", "0.001699903579759509", "0.001699903579759509", "2", "0.6666666666666666", "1.734127333333333", "5.202382", "0.07488702167465976", "0.07488702167465977", "12.17086228289311", "109.537760546038", "0.000127300716222282", "0.0001273007162222821"),
    @("ECs", "Postn", "Ptk7", "FAPs", "3", "1", "7.018436333333334", "21.055309", "0.001699903579759509", "0.001699903579759509", "3", "1", "11.370088", "34.110264", "0.4910089415764484", "0.4910089415764484", "79.80023873239735", "718.202148591576", "0.0008346678574797321", "0.0008346678574797323"),
    @("ECs", "Postn", "Ptk7", "M2", "3", "1", "7.018436333333334", "21.055309", "0.001699903579759509", "0.001699903579759509", "2", "0.6666666666666666", "0.02674833333333333", "0.080245", "0.001155107228627785", "0.001155107228627785", "0.1877314745227778", "1.689583270705", "1.963570912950457E-06", "1.963570912950457E-06"),
    @("ECs", "Postn", "Ptk7", "sCs", "3", "1", "7.018436333333334", "21.055309", "0.001699903579759509", "0.001699903579759509", "3", "1", "10.02561666666667", "30.07685", "0.432948929520264", "0.4329489295202641", "70.36415227740555", "633.2773704966501", "0.0007359714351445441", "0.0007359714351445443"),
    @("FAPs", "Postn", "Ptk7", "ECs", "3", "1", "3983.399495666667", "11950.198487", "0.9648010953858698", "0.96480109538587", "2", "0.6666666666666666", "1.734127333333333", "5.202382", "0.07488702167465976", "0.07488702167465977", "6907.721945021782", "62169.49750519604", "0.07225108054189711", "0.07225108054189713"),
    @("FAPs", "Postn", "Ptk7", "FAPs", "3", "1", "3983.399495666667", "11950.198487", "0.9648010953858698", "0.96480109538587", "3", "1", "11.370088", "34.110264", "0.4910089415764484", "0.4910089415764484", "45291.60280488562", "407624.4252439706", "0.4737259646772139", "0.473725964677214"),
    @("FAPs", "Postn", "Ptk7", "M2", "3", "1", "3983.399495666667", "11950.198487", "0.9648010953858698", "0.96480109538587", "2", "0.6666666666666666", "0.02674833333333333", "0.080245", "0.001155107228627785", "0.001155107228627785", "106.5492975099239", "958.9436775893149", "0.001114448719468223", "0.001114448719468224"),
    @("FAPs", "Postn", "Ptk7", "sCs", "3", "1", "3983.399495666667", "11950.198487", "0.9648010953858698", "0.96480109538587", "3", "1", "10.02561666666667", "30.07685", "0.432948929520264", "0.4329489295202641", "39936.03637374732", "359424.3273637259", "0.4177096014472905", "0.4177096014472906"),
    @("M2", "Postn", "Ptk7", "ECs", "3", "1", "1.256405", "3.769215", "0.0003043081472413079", "0.000304308147241308", "2", "0.6666666666666666", "1.734127333333333", "5.202382", "0.07488702167465976", "0.07488702167465977", "2.178766252236667", "19.60889627013", "2.278873081823538E-05", "2.278873081823539E-05"),
    @("M2", "Postn", "Ptk7", "FAPs", "3", "1", "1.256405", "3.769215", "0.0003043081472413079", "0.000304308147241308", "3", "1", "11.370088", "34.110264", "0.4910089415764484", "0.4910089415764484", "14.28543541364", "128.56891872276", "0.0001494180212900446", "0.0001494180212900446"),
    @("M2", "Postn", "Ptk7", "M2", "3", "1", "1.256405", "3.769215", "0.0003043081472413079", "0.000304308147241308", "2", "0.6666666666666666", "0.02674833333333333", "0.080245", "0.001155107228627785", "0.001155107228627785", "0.03360673974166666", "0.302460657675", "3.515085406087631E-07", "3.515085406087632E-07"),
    @("M2", "Postn", "Ptk7", "sCs", "3", "1", "1.256405", "3.769215", "0.0003043081472413079", "0.000304308147241308", "3", "1", "10.02561666666667", "30.07685", "0.432948929520264", "0.4329489295202641", "12.59623490808333", "113.36611417275", "0.0001317498865924192", "0.0001317498865924192"),
    @("sCs", "Postn", "Ptk7", "ECs", "3", "1", "137.051796", "411.155388", "0.03319469288712931", "0.03319469288712931", "2", "0.6666666666666666", "1.734127333333333", "5.202382", "0.07488702167465976", "0.07488702167465977", "237.665265526024", "2138.987389734216", "0.002485851685722127", "0.002485851685722127"),
    @("sCs", "Postn", "Ptk7", "FAPs", "3", "1", "137.051796", "411.155388", "0.03319469288712931", "0.03319469288712931", "3", "1", "11.370088", "34.110264", "0.4910089415764484", "0.4910089415764484", "1558.290981078048", "14024.61882970243", "0.01629889102046462", "0.01629889102046462"),
    @("sCs", "Postn", "Ptk7", "M2", "3", "1", "137.051796", "411.155388", "0.03319469288712931", "0.03319469288712931", "2", "0.6666666666666666", "0.02674833333333333", "0.080245", "0.001155107228627785", "0.001155107228627785", "3.66590712334", "32.99316411006", "3.834342970600237E-05", "3.834342970600238E-05"),
    @("sCs", "Postn", "Ptk7", "sCs", "3", "1", "137.051796", "411.155388", "0.03319469288712931", "0.03319469288712931", "3", "1", "10.02561666666667", "30.07685", "0.432948929520264", "0.4329489295202641", "1374.0287701742", "12366.2589315678", "0.01437160675123656", "0.01437160675123656")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $col = $c + 1
        $cell = $ws.Cells.Item($r, $col)
        if ($col -le 4) {
            $cell.Value = $rowValues[$c]
        } else {
            $cell.Value = [double]$rowValues[$c]
        }
    }
}

